$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 177937
$ws.Range("C4").Value = 167893
$ws.Range("C7").Value = 5.64
$ws.Range("C8").Value = 65
